# Switch license from BY-NC to BY-SA
#
# The credits/license slide (slide 2, "Title 1" placeholder) states the
# deck is licensed "CC BY-NC 4.0" and links to the corresponding
# creativecommons.org URL. Update both the visible license label and the
# visible hyperlink text to the BY-SA variant, leaving everything else
# (including the hyperlink's underlying target relationship) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# "CC BY-NC 4.0. To view a copy of this license, visit " -> "CC BY-SA 4.0. ..."
$full = $tr.Text
$idx = $full.IndexOf("BY-NC ")
$sub = $tr.Characters($idx + 1, 6)
$sub.Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0" -> ".../by-sa/4.0"
$full = $tr.Text
$idx = $full.IndexOf("creativecommons.org/licenses/by-nc/4.0")
$sub = $tr.Characters($idx + 1, 39)
$sub.Text = "creativecommons.org/licenses/by-sa/4.0"
